# Update the first-page header (header3.xml) signature block:
# "Dr. med. Thiên-Trí Lâm" -> "PD Dr. med. Thiên-Trí Lâm"

$d = $word.ActiveDocument

$sections = $d.Sections
for ($i = 1; $i -le $sections.Count; $i++) {
    $section = $sections.Item($i)
    $header = $section.Headers.Item(1)  # wdHeaderFooterPrimary = 1 (first-page header when titlePg set)
    if ($header.Exists) {
        $header.Range.Find.Execute("Dr. med. Thiên-Trí Lâm", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "PD Dr. med. Thiên-Trí Lâm", 2)
    }
}
